# "modify doc for 0416" - populate the (previously blank) 3rd slide
# ("进度规划" / progress-planning slide) with its title and body content.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# ---------------------------------------------------------------------
# Title placeholder ("标题 1", shape 1): "进度规划"
# ---------------------------------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "进度规划"
$title.LanguageID = "zh-CN"

# ---------------------------------------------------------------------
# Content placeholder ("内容占位符 2", shape 2): two paragraphs of notes
# followed by two blank paragraphs.
# ---------------------------------------------------------------------
$body = $s.Shapes.Item(2).TextFrame.TextRange

$para1Text = "在ARM secure world 移植zCore（下周内要完成，可以多问杨德睿学长和贾越凯学长）"
$para2Text = "通过 shared memory 设计实现REE/TEE的IPC，自行设计对应的demo和测例（也需要多问助教和学长）"

# Build the 4 paragraphs: two content paragraphs + two trailing blank ones.
$body.Text = $para1Text + "`r" + $para2Text + "`r`r"

# --- paragraph 1 run formatting -----------------------------------------
$para1 = $body.Paragraphs(1, 1)

$runsP1 = @(
    @{ Text = "在"; Lang = "zh-CN" },
    @{ Text = "ARM secure world "; Lang = "en-US" },
    @{ Text = "移植"; Lang = "zh-CN" },
    @{ Text = "zCore"; Lang = "en-US" },
    @{ Text = "（下周内要完成，可以多问杨德睿学长和贾越凯学长）"; Lang = "zh-CN" }
)

$pos = 1
foreach ($run in $runsP1) {
    $len = $run.Text.Length
    $rng = $para1.Characters($pos, $len)
    $rng.LanguageID = $run.Lang
    $rng.Font.Bold = 0
    $rng.Font.Italic = 0
    $rng.Font.Shadow = 0
    $rng.Font.Name = "-apple-system"
    $pos += $len
}

# endParaRPr for paragraph 1 (after "）") picks up en-US / b=0 i=0 styling;
# touch the trailing (zero-length) position so the paragraph-mark run
# properties match.
$end1 = $para1.Characters($para1.Length + 1, 0)
$end1.LanguageID = "en-US"
$end1.Font.Bold = 0
$end1.Font.Italic = 0
$end1.Font.Shadow = 0
$end1.Font.Name = "-apple-system"

# --- paragraph 2 run formatting -----------------------------------------
$para2 = $body.Paragraphs(2, 1)

$runsP2 = @(
    @{ Text = "通过 "; Lang = "zh-CN"; Dim = $true },
    @{ Text = "shared memory "; Lang = "en-US"; Dim = $true },
    @{ Text = "设计实现"; Lang = "zh-CN"; Dim = $true },
    @{ Text = "REE/TEE"; Lang = "en-US"; Dim = $true },
    @{ Text = "的"; Lang = "zh-CN"; Dim = $true },
    @{ Text = "IPC"; Lang = "en-US"; Dim = $true },
    @{ Text = "，"; Lang = "zh-CN"; Dim = $true },
    @{ Text = "自行设计对应的"; Lang = "zh-CN"; Dim = $false },
    @{ Text = "demo"; Lang = "en-US"; Dim = $false },
    @{ Text = "和测例（也需要多问助教和学长）"; Lang = "zh-CN"; Dim = $false }
)

$pos = 1
foreach ($run in $runsP2) {
    $len = $run.Text.Length
    $rng = $para2.Characters($pos, $len)
    $rng.LanguageID = $run.Lang
    if ($run.Dim) {
        $rng.Font.Bold = 0
        $rng.Font.Italic = 0
    }
    $rng.Font.Shadow = 0
    $rng.Font.Name = "-apple-system"
    $pos += $len
}

# endParaRPr for paragraph 2 matches the "dim" (b=0 i=0) en-US styling.
$end2 = $para2.Characters($para2.Length + 1, 0)
$end2.LanguageID = "en-US"
$end2.Font.Bold = 0
$end2.Font.Italic = 0
$end2.Font.Shadow = 0
$end2.Font.Name = "-apple-system"

# --- paragraph 3 (blank) --------------------------------------------------
$para3 = $body.Paragraphs(3, 1)
$end3 = $para3.Characters(1, 0)
$end3.LanguageID = "en-US"
$end3.Font.Bold = 0
$end3.Font.Italic = 0
$end3.Font.Shadow = 0
$end3.Font.Name = "-apple-system"

# paragraph 4 (blank, zh-CN) is left as the pre-existing trailing empty
# paragraph and needs no further changes.
